$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.493305087089539
$ws.Range("B1").Value = 1.559166789054871
$ws.Range("C1").Value = 1.749374032020569
$ws.Range("D1").Value = 2.723487854003906
$ws.Range("E1").Value = 4.317018032073975
